$d = $word.ActiveDocument

# 1. Paragraph 6: Objetivos content
$p = $d.Paragraphs.Item(6).Range
$dup = $p.Duplicate
$found = $dup.Find.Execute("A disciplina busca introduzir o aluno ao ambiente de engenharia, propondo problemas desafiadores gerando aptidão para solução de problemas. Apresentar a Engenharia de Materiais e seus campos de atuação, aspectos legais e éticos, bem como o mercado de trabalho para o engenheiro de materiais no Século XXI. Propiciar aos alunos uma visão geral do curso, com apresentação do currículo do curso de Engenharia de Materiais da EEL. Apresentar aos alunos uma visão da evolução histórica dos materiais com o homem. Descrever exemplos marcantes da introdução de novos materiais e as mudanças sociais provocadas. Apresentar o caráter interdisciplinar da Ciência e Engenharia de Materiais e suas ligações com outros ramos da Ciência. Apresentar estudos de caso demonstrando este caráter interdisciplinar.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { Write-Host "WARNING: locate 1 (Objetivos content) in paragraph 6 not found!" }
$target = $d.Range($dup.Start, $dup.End)
$target.Text = "1- A importância dos materiais na evolução do homem na pré-história. Alquimia, Revolução Científica e a Revolução Industrial. 2-O Engenheiro como um profissional, funções da engenharia, a ética e comunicação na engenharia 3-A grandes áreas da Engenharia de Materiais. A interdisciplinaridade da Ciência e Engenharia de Materiais. 4- Perspectivas para a Engenharia de Materiais no século XXI. 5- O currículo do curso de engenharia de materiais da EEL-USP. 6- Noções básicas de Projetos em Engenharia.``vEm todos o conteúdo do curso serão abordados aspectos sociais, ambientais, éticos, legais e econômicos para ampliar as competências dos alunos"

# 2. Paragraph 8: Docentes run1
$p = $d.Paragraphs.Item(8).Range
$dup = $p.Duplicate
$found = $dup.Find.Execute("984972 - Hugo Ricardo Zschommler Sandim", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { Write-Host "WARNING: locate 2 (Docentes run1) in paragraph 8 not found!" }
$target = $d.Range($dup.Start, $dup.End)
$target.Text = "A disciplina busca introduzir o aluno ao ambiente de engenharia, propondo problemas desafiadores gerando aptidão para solução de problemas. Apresentar a Engenharia de Materiais e seus campos de atuação, aspectos legais e éticos, bem como o mercado de trabalho para o engenheiro de materiais no Século XXI. Propiciar aos alunos uma visão geral do curso, com apresentação do currículo do curso de Engenharia de Materiais da EEL. Apresentar aos alunos uma visão da evolução histórica dos materiais com o homem. Descrever exemplos marcantes da introdução de novos materiais e as mudanças sociais provocadas. Apresentar o caráter interdisciplinar da Ciência e Engenharia de Materiais e suas ligações com outros ramos da Ciência. Apresentar estudos de caso demonstrando este caráter interdisciplinar."

# 3. Paragraph 8: Docentes run2
$p = $d.Paragraphs.Item(8).Range
$dup = $p.Duplicate
$found = $dup.Find.Execute("7459752 - Maria Ismenia Sodero Toledo Faria", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { Write-Host "WARNING: locate 3 (Docentes run2) in paragraph 8 not found!" }
$target = $d.Range($dup.Start, $dup.End)
$target.Text = "As características importantes de um engenheiro: aptidões interpessoais, aptidões de comunicação, liderança e competência. O engenheiro, profissional que busca solucionar problemas. 2-A Engenharia de Materiais: áreas de atuação e mercado de trabalho. Aplicação. A importância dos materiais na evolução do homem, as grandes áreas e interdisciplinaridade da Ciência e Engenharia de Materiais. Visita ao Departamento de Engenharia de Materiais. Conhecimento dos Grupos de Pesquisa do Departamento. Perspectivas para a Engenharia de Materiais no século XXI. 3- O campo de trabalho do engenheiro de materiais e suas áreas de atuação. Visita externa para integralização dos conhecimentos. 4- O currículo do curso de engenharia de materiais na EEL/USP. 5- Apresentação do método de trabalho com projetos, definindo os atributos de um projeto de engenharia, mapas conceituais e ferramentas que ilustram ideias e relações entre elas. Formular estratégias para resolução de problemas de engenharia. Estudo de casos. Viagem didática complementar."

# 4. Paragraph 10: ProgramaResumido content
$p = $d.Paragraphs.Item(10).Range
$dup = $p.Duplicate
$found = $dup.Find.Execute("1- A importância dos materiais na evolução do homem na pré-história. Alquimia, Revolução Científica e a Revolução Industrial. 2-O Engenheiro como um profissional, funções da engenharia, a ética e comunicação na engenharia 3-A grandes áreas da Engenharia de Materiais. A interdisciplinaridade da Ciência e Engenharia de Materiais. 4- Perspectivas para a Engenharia de Materiais no século XXI. 5- O currículo do curso de engenharia de materiais da EEL-USP. 6- Noções básicas de Projetos em Engenharia.^lEm todos o conteúdo do curso serão abordados aspectos sociais, ambientais, éticos, legais e econômicos para ampliar as competências dos alunos", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { Write-Host "WARNING: locate 4 (ProgramaResumido content) in paragraph 10 not found!" }
$target = $d.Range($dup.Start, $dup.End)
$target.Text = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"

# 5. Paragraph 12: Programa content
$p = $d.Paragraphs.Item(12).Range
$dup = $p.Duplicate
$found = $dup.Find.Execute("As características importantes de um engenheiro: aptidões interpessoais, aptidões de comunicação, liderança e competência. O engenheiro, profissional que busca solucionar problemas. 2-A Engenharia de Materiais: áreas de atuação e mercado de trabalho. Aplicação. A importância dos materiais na evolução do homem, as grandes áreas e interdisciplinaridade da Ciência e Engenharia de Materiais. Visita ao Departamento de Engenharia de Materiais. Conhecimento dos Grupos de Pesquisa do Departamento. Perspectivas para a Engenharia de Materiais no século XXI. 3- O campo de trabalho do engenheiro de materiais e suas áreas de atuação. Visita externa para integralização dos conhecimentos. 4- O currículo do curso de engenharia de materiais na EEL/USP. 5- Apresentação do método de trabalho com projetos, definindo os atributos de um projeto de engenharia, mapas conceituais e ferramentas que ilustram ideias e relações entre elas. Formular estratégias para resolução de problemas de engenharia. Estudo de casos. Viagem didática complementar.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { Write-Host "WARNING: locate 5 (Programa content) in paragraph 12 not found!" }
$target = $d.Range($dup.Start, $dup.End)
$target.Text = "Média Aritmética dos Projetos, Trabalhos, Relatórios e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."

# 6. Paragraph 14: Avaliacao/Metodo answer
$p = $d.Paragraphs.Item(14).Range
$dup = $p.Duplicate
$found = $dup.Find.Execute("Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { Write-Host "WARNING: locate 6 (Avaliacao/Metodo answer) in paragraph 14 not found!" }
$target = $d.Range($dup.Start, $dup.End)
$target.Text = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

# 7. Paragraph 14: Avaliacao/Criterio answer
$p = $d.Paragraphs.Item(14).Range
$dup = $p.Duplicate
$found = $dup.Find.Execute("Média Aritmética dos Projetos, Trabalhos, Relatórios e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { Write-Host "WARNING: locate 7 (Avaliacao/Criterio answer) in paragraph 14 not found!" }
$target = $d.Range($dup.Start, $dup.End)
$target.Text = "1) BROCKMAN, J.B. Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2010.``v2) M.T. HOLTZAPPLE, W.D. REECE, Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2006.``v2) CALLISTER Jr., W.D. Ciência e Engenharia de Materiais: Uma Introdução. LTC Livros Científicos Editora, 7a.ed., 2008. ``v4) - COHEN, M. (Ed.). Ciência e Engenharia de Materiais: sua Evolução, Prática e Perspectivas. Parte I: Materiais na história e na sociedade, 98p. Parte II: A Ciência e Engenharia de Materiais como uma multidisciplina, Tradução: José Roberto Gonçalves da Silva, São Carlos, UFSCar, 1985.``v5) Artigos científicos"

# 8. Paragraph 14: Avaliacao/Norma answer
$p = $d.Paragraphs.Item(14).Range
$dup = $p.Duplicate
$found = $dup.Find.Execute("NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { Write-Host "WARNING: locate 8 (Avaliacao/Norma answer) in paragraph 14 not found!" }
$target = $d.Range($dup.Start, $dup.End)
$target.Text = "984972 - Hugo Ricardo Zschommler Sandim"

# 9. Paragraph 16: Bibliografia content
$p = $d.Paragraphs.Item(16).Range
$dup = $p.Duplicate
$found = $dup.Find.Execute("1) BROCKMAN, J.B. Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2010.^l2) M.T. HOLTZAPPLE, W.D. REECE, Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2006.^l2) CALLISTER Jr., W.D. Ciência e Engenharia de Materiais: Uma Introdução. LTC Livros Científicos Editora, 7a.ed., 2008. ^l4) - COHEN, M. (Ed.). Ciência e Engenharia de Materiais: sua Evolução, Prática e Perspectivas. Parte I: Materiais na história e na sociedade, 98p. Parte II: A Ciência e Engenharia de Materiais como uma multidisciplina, Tradução: José Roberto Gonçalves da Silva, São Carlos, UFSCar, 1985.^l5) Artigos científicos", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found) { Write-Host "WARNING: locate 9 (Bibliografia content) in paragraph 16 not found!" }
$target = $d.Range($dup.Start, $dup.End)
$target.Text = "7459752 - Maria Ismenia Sodero Toledo Faria"
